$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that had their values removed (column R period 31/12/2002 dropped
# for the income-statement rows, plus the minority-interest row 78 cleared of its zero placeholders).
$clearCells = @("R58","R59","R60","R61","R62","R63","R64","R65","R66","R67","R68","R69","R70","R71","R72","R73","R74","R75","R76","R77","B78","D78","E78","F78","G78","H78","I78","J78","K78","L78","M78","N78","O78","P78","Q78","R78","S78","T78","U78","V78","W78","X78","Y78","Z78","AA78","AB78","AC78","AD78","AE78","AF78","AG78","AH78","AI78","R79")
foreach ($addr in $clearCells) {
    $ws.Range($addr).ClearContents()
}

# Apply the small recalculated value corrections produced by the concatenation merge.
$valueUpdates = @{
    "O58" = 777193.856
    "Z58" = 683779.968
    "O59" = -295961.152
    "AD59" = -413632.032
    "AH59" = -332652.064
    "AD60" = 625486.0159999999
    "K61" = -224821.968
    "V63" = -265701.008
    "AH63" = -165348
    "O64" = -156670.992
    "AH64" = -124617.992
    "V66" = -131853
    "AD66" = 113165.992
    "V67" = 47574
    "K68" = 152038
    "V68" = -184078
    "AH69" = 120343.992
    "K73" = 95462
    "Z73" = 126655.992
    "Z79" = 13905
}
foreach ($addr in $valueUpdates.Keys) {
    $ws.Range($addr).Value = $valueUpdates[$addr]
}
